$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 15 (the previous file had a gap between
# row 14 and row 16, so row 15 did not exist yet). This shifts the old
# rows 16..23 down to 17..24 and copies formatting from row 14 above,
# matching the surrounding rows' style.
$ws.Rows.Item(15).Insert()

# Row-insert also materializes column D's default style (style 5) into the
# new row even though it should stay empty - clear it so no D15 cell remains.
$ws.Cells.Item(15, 4).Clear()

# Fill in the new "devices" hash documentation row.
$ws.Cells.Item(15, 2).Value = "gbridge:u{userid}:devices"
# Leading apostrophe forces Excel's "stored as text" quote-prefix formatting,
# since the text begins with three dashes.
$ws.Cells.Item(15, 3).Value = "'--- (simple SET/GET)"
$ws.Cells.Item(15, 6).Value = "Array containing information about user's devices and their supported traits"

# Keep the merged "gbridge:u{userid}:d0:grequest" column's centered/middle
# alignment consistent across the (now taller) merged block.
$ws.Range("B10:B14").HorizontalAlignment = -4108
$ws.Range("B10:B14").VerticalAlignment = -4108
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("B15").VerticalAlignment = -4108

# Match the author's final selection/active cell.
$ws.Range("F15").Select()
